$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===================================================================
# PART 1 - Table 2 (rows 20-34): fill in the previously-blank columns
# I (HEX2DEC of the G column) and J (duty-cycle-period as a decimal)
# ===================================================================

# Fix up cell styles first so they match the final look:
#  - I column should be style s=3 all the way down (I25:I29 were s=5)
#  - J column should be style s=5 all the way down (J20:J24 & J30:J34 were s=3)
$ws.Range("I20").Copy()
$ws.Range("I25:I29").PasteSpecial(-4122)

$ws.Range("J25").Copy()
$ws.Range("J20:J24").PasteSpecial(-4122)
$ws.Range("J30:J34").PasteSpecial(-4122)

# Now fill in the formulas (letting Excel relative-adjust + auto-share them)
$ws.Range("I20").Formula = "=HEX2DEC(G20)"
$ws.Range("I21:I34").Formula = "=HEX2DEC(G21)"

$ws.Range("J20:J28").Formula = "=(B20*E20)/((1/C20)*16)"
$ws.Range("J29").Formula = "=(B29*E29)/((1/C29)*16)"
$ws.Range("J30:J34").Formula = "=(B30*E30)/((1/C30)*16)"

Write-Host "Table 2 I/J columns filled in."

# ===================================================================
# PART 2 - Table 3 (rows 36-51): a brand new table, structurally a
# copy of table 2 (FREQUENCIES / DUTY CYCLE / FOSC / PRESCALER /
# PERIOD / PR2 / CCPR1L-hex / duty-cycle-period), plus two extra
# columns: I = HEX2DEC(G) and J = comparison flag against table 2.
# ===================================================================

# --- header row 36: copy formatting from row 19's header, then set text
$headerCols = @("A","B","C","D","E","F","G","H","I")
foreach ($c in $headerCols) {
    $ws.Range("$c" + "19").Copy()
    $ws.Range("$c" + "36").PasteSpecial(-4122)
}
$ws.Range("A36").Value = "FREQUENCIES"
$ws.Range("B36").Value = "DUTY CYCLE"
$ws.Range("C36").Value = "FOSC"
$ws.Range("D36").Value = "PRESCALER"
$ws.Range("E36").Value = "PERIOD"
$ws.Range("F36").Value = "PR2"
$ws.Range("G36").Value = "(CCPR1L : CCP1CON < 5 : 4 > (HEX)"
$ws.Range("H36").Value = "DUTY CYCLE PERIOD"

Write-Host "Table 3 header row done."

# --- body rows 37-51: three blocks of 5 rows (300Hz, 500Hz, 1000Hz),
# each block copying the column formatting from its table-2 counterpart.
$blocks = @(
    @{ Freq = 300;  Dest = 37; FmtSrc = 20 },
    @{ Freq = 500;  Dest = 42; FmtSrc = 25 },
    @{ Freq = 1000; Dest = 47; FmtSrc = 30 }
)
$duties = @(0.1, 0.25, 0.5, 0.75, 0.95)
$bodyCols = @("A","B","C","D","E","F","G","H")

foreach ($block in $blocks) {
    $destStart = $block.Dest
    $destEnd = $destStart + 4
    $fmtSrc = $block.FmtSrc

    # copy the number formatting (fill/font/alignment) for A:H from the
    # matching table-2 row block down onto the new rows
    foreach ($c in $bodyCols) {
        $ws.Range("$c" + "$fmtSrc").Copy()
        $ws.Range("$c" + "$destStart" + ":" + "$c" + "$destEnd").PasteSpecial(-4122)
    }
    # I column formatting always matches I20 (style s=3)
    $ws.Range("I20").Copy()
    $ws.Range("I$destStart" + ":" + "I$destEnd").PasteSpecial(-4122)

    # values for A (freq), B (duty cycle), C (fosc), D (prescaler)
    $row = $destStart
    foreach ($d in $duties) {
        $ws.Range("A$row").Value = $block.Freq
        $ws.Range("B$row").Value = $d
        $ws.Range("C$row").Value = 4000000
        $ws.Range("D$row").Value = 16
        $row++
    }

    # formulas for E, F, H (same shape as table 2's corresponding columns)
    $ws.Range("E$destStart").Formula = "=1/A$destStart"
    if ($destEnd -gt $destStart) {
        $ws.Range("E" + ($destStart + 1) + ":E$destEnd").Formula = "=1/A" + ($destStart + 1)
    }

    $ws.Range("F$destStart" + ":F$destEnd").Formula = "=DEC2HEX(ROUND((E$destStart/(4*(2.5*10^-7)*D$destStart))-1,0))"

    $ws.Range("H$destStart" + ":H$destEnd").Formula = "=(B$destStart*E$destStart)"

    # G column: deliberately NOT rounded (matches the author's formula)
    $ws.Range("G$destStart" + ":G$destEnd").Formula = "=DEC2HEX((B$destStart*E$destStart)/((1/C$destStart)*16))"

    # I column: HEX2DEC of this row's G cell
    $ws.Range("I$destStart" + ":I$destEnd").Formula = "=HEX2DEC(G$destStart)"

    # J column: compares this row's G value against the corresponding
    # table-2 row (17 rows above) - no explicit style (default)
    $refRow = $destStart - 17
    $ws.Range("J$destStart" + ":J$destEnd").Formula = "=IF(G$refRow=G$destStart,1,0)"
}

Write-Host "Table 3 body rows done."

# ===================================================================
# PART 3 - cosmetic: move the active selection close to where the
# author left it (best-effort; window geometry itself isn't scriptable)
# ===================================================================
$ws.Range("H17").Select()

Write-Host "Done."
